# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Both sheets share the same row layout; only row 26's new values differ
# from the old values (4169 -> 4184 on 展览, 4170 -> 4184 on 全部类型), but
# since we're only setting the target/new values, we just apply the same
# new value to both sheets for every row.

$wb = $excel.ActiveWorkbook

# row -> new F-column value
$updates = @{
    2  = 136
    4  = 62
    6  = 128
    7  = 1263
    8  = 1535
    10 = 390
    12 = 148
    17 = 302
    18 = 323
    19 = 1734
    23 = 667
    25 = 335
    26 = 4184
    28 = 269
    29 = 1086
    30 = 486
    32 = 540
    34 = 248
    36 = 137
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
